$d = $word.ActiveDocument

# Bump every existing paragraph's run(s) and paragraph mark up to 14pt
# (sz=28 half-points) / 18pt complex-script (szCs=36 half-points).
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    $p.Range.Font.Size = 14
    $p.Range.Font.SizeBi = 18
}

# The "_GoBack" bookmark currently sits at the end of the last paragraph
# ("For Batch Size 32 in VGG16, Exceeded GPU memory."). It needs to move
# to the new paragraph being appended below, so drop it from here first.
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

# Append a new list item after the last paragraph for the new problem
# faced, inheriting the same list/numbering + rPr formatting.
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs($d.Paragraphs.Count)
$newPara.Range.Text = "Cropping images of different sizes from different sources was an issue."

# Re-create the "_GoBack" bookmark, collapsed at the start of the new
# paragraph (before its run), matching the original bookmarkStart/End
# placement ahead of the text run.
$bmRange = $newPara.Range.Duplicate
$bmRange.Collapse(1)
$d.Bookmarks.Add("_GoBack", $bmRange)
